$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "26.148.15"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "1.668.46"
$ws.Range("E3").Value = "  -0.81%  "
Set-TextValue "D4" "1.002"
$ws.Range("E4").Value = "  -0.49%  "
Set-TextValue "D5" "209.64"
$ws.Range("E5").Value = "  -2.92%  "
Set-TextValue "D6" "0.5236"
$ws.Range("E6").Value = "  -0.51%  "
Set-TextValue "D7" "1.003"
$ws.Range("E7").Value = "  -0.46%  "
Set-TextValue "D8" "0.2623"
$ws.Range("E8").Value = "  -2.72%  "
Set-TextValue "D9" "0.06331"
$ws.Range("E9").Value = "  -0.65%  "
$ws.Range("E10").Value = "  -1.66%  "
Set-TextValue "D11" "0.07534"
$ws.Range("E11").Value = "  -1.49%  "
$ws.Range("D12").Value = "1.672.49"
$ws.Range("E12").Value = "  -0.44%  "
$ws.Range("E13").Value = "  -1.44%  "
$ws.Range("E14").Value = "  -4.48%  "
Set-TextValue "D15" "66.40"
$ws.Range("E15").Value = "  +0.52%  "
Set-TextValue "D16" "0.000007962"
$ws.Range("E16").Value = "  -4.44%  "
$ws.Range("D17").Value = "26.152.07"
$ws.Range("E17").Value = "  -0.39%  "
$ws.Range("E18").Value = "  -0.51%  "
Set-TextValue "D19" "4.755"
$ws.Range("E19").Value = "  -2.23%  "
Set-TextValue "D20" "186.54"
$ws.Range("E20").Value = "  -1.35%  "
Set-TextValue "D21" "10.30"
$ws.Range("E21").Value = "  -4.70%  "
Set-TextValue "D22" "6.190"
$ws.Range("E22").Value = "  -0.60%  "
$ws.Range("E23").Value = "  -0.56%  "
Set-TextValue "D24" "149.28"
$ws.Range("E24").Value = "  +0.35%  "
Set-TextValue "D25" "0.1249"
$ws.Range("E25").Value = "  -0.88%  "
$ws.Range("E26").Value = "  -3.68%  "
Set-TextValue "D27" "15.85"
$ws.Range("E27").Value = "  +0.84%  "
Set-TextValue "D28" "0.06378"
$ws.Range("E28").Value = "  +1.55%  "
Set-TextValue "D29" "1.349"
$ws.Range("E29").Value = "  -1.68%  "
$ws.Range("E30").Value = "  -3.13%  "
Set-TextValue "D31" "3.519"
$ws.Range("E31").Value = "  -1.34%  "
Set-TextValue "D32" "3.412"
$ws.Range("E32").Value = "  -4.13%  "
Set-TextValue "D33" "1.648"
$ws.Range("E33").Value = "  -2.06%  "
Set-TextValue "D34" "1.005"
$ws.Range("E34").Value = "  -1.63%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D35" "2.409"
$ws.Range("E35").Value = "  -0.41%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D36" "0.6010"
$ws.Range("E36").Value = "  -1.83%  "
Set-TextValue "D37" "2.739"
$ws.Range("E37").Value = "  -0.54%  "
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "1.114.50"
$ws.Range("E38").Value = "  +1.40%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D39" "6.151"
$ws.Range("E39").Value = "  -0.42%  "
Set-TextValue "D40" "0.01616"
$ws.Range("E40").Value = "  -0.19%  "
Set-TextValue "D41" "0.8656"
$ws.Range("E41").Value = "  -2.96%  "
Set-TextValue "D42" "1.004"
$ws.Range("E42").Value = "  -0.70%  "
Set-TextValue "D43" "100.37"
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("D44").Value = "1.824.02"
$ws.Range("E44").Value = "  -0.49%  "
Set-TextValue "D45" "0.00000000110"
$ws.Range("E45").Value = "  +0.32%  "
Set-TextValue "D46" "55.52"
$ws.Range("E46").Value = "  -3.19%  "
$ws.Range("E47").Value = "  -0.01%  "
Set-TextValue "D48" "8.072"
$ws.Range("E48").Value = "  +0.06%  "
Set-TextValue "D49" "0.05232"
$ws.Range("E49").Value = "  -0.85%  "
Set-TextValue "D50" "0.4238"
$ws.Range("E50").Value = "  -1.15%  "
Set-TextValue "D51" "5.918"
$ws.Range("E51").Value = "  -1.57%  "
